$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 ---
$ws.Range("A24").Value = 111958182
$ws.Range("B24").Value = 55611
$ws.Range("C24").Value = "Ovaliderad"
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 102612
$ws.Range("F24").Value = "Järpe"
$ws.Range("G24").Value = "Tetrastes bonasia"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("I24").NumberFormat = "@"
$ws.Range("I24").Value = "3"
$ws.Range("K24").Value = ""
$ws.Range("L24").Value = "hona"
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = ""
$ws.Range("P24").Value = "Österåsen, Österås, Ång"
$ws.Range("Q24").Value = 609746.731343443
$ws.Range("R24").Value = 7011953.229753771
$ws.Range("S24").Value = 10
$ws.Range("T24").Value = "Västernorrland"
$ws.Range("U24").Value = "Sollefteå"
$ws.Range("V24").Value = "Ångermanland"
$ws.Range("W24").Value = "Ed"
$ws.Range("Y24").Value = "2023-09-04"
$ws.Range("Z24").Value = "00:00"
$ws.Range("AA24").Value = "2023-09-04"
$ws.Range("AB24").Value = "00:00"
$ws.Range("AC24").Value = "1K"
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AT24").Value = ""
$ws.Range("AW24").Value = "Lennart Vessberg"
$ws.Range("AX24").Value = "Lennart Vessberg"
$ws.Range("AY24").Value = ""

# --- Row 25 ---
$ws.Range("A25").Value = 111958205
$ws.Range("B25").Value = 96348
$ws.Range("C25").Value = "Ovaliderad"
$ws.Range("D25").Value = "VU"
$ws.Range("E25").Value = 220787
$ws.Range("F25").Value = "Knärot"
$ws.Range("G25").Value = "Goodyera repens"
$ws.Range("H25").Value = "(L.) R. Br."
$ws.Range("I25").NumberFormat = "@"
$ws.Range("I25").Value = "1"
$ws.Range("J25").Value = "plantor/tuvor"
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("P25").Value = "Österåsen, Ång"
$ws.Range("Q25").Value = 609802.6803741428
$ws.Range("R25").Value = 7011969.124995505
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = "Västernorrland"
$ws.Range("U25").Value = "Sollefteå"
$ws.Range("V25").Value = "Ångermanland"
$ws.Range("W25").Value = "Ed"
$ws.Range("Y25").Value = "2023-09-04"
$ws.Range("Z25").Value = "00:00"
$ws.Range("AA25").Value = "2023-09-04"
$ws.Range("AB25").Value = "00:00"
$ws.Range("AC25").Value = "½ m2"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AF25").Value = ""
$ws.Range("AG25").Value = $false
$ws.Range("AT25").Value = ""
$ws.Range("AW25").Value = "Lennart Vessberg"
$ws.Range("AX25").Value = "Lennart Vessberg"
$ws.Range("AY25").Value = ""
